# ============================================================
# Edit script: insert a new "2022-Q3" sheet with fund holding
# data (after "总计"), and prepend a corresponding summary row
# on the "总计" (totals) sheet, renumbering its index column.
# ============================================================

function Set-TextCell($ws, $row, $col, $val) {
    # Force the value to be stored as text (quote-prefixed) even
    # when it looks like a number (e.g. "011296" or "93.82"),
    # matching the source workbook's string-cell representation
    # and preventing loss of leading/trailing zeros.
    $ws.Cells.Item($row, $col).Value = "'" + $val
}

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
# A sheet that already has the fund-table header/index formatting
# we want to reuse (the old "2022-Q2" sheet, currently 2nd sheet).
$formatSrc = $wb.Worksheets.Item(2)

# ------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Header row
$q3Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q3Headers.Count; $c++) {
    $q3Sheet.Cells.Item(1, $c + 2).Value = $q3Headers[$c]
}

# Data rows: index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @(0, "450004", "国富深化价值混合", "93.82", "84.25", "2.39", "2.2423", 4),
    @(1, "213001", "宝盈鸿利收益灵活配置混合A", "16.06", "87.42", "5.17", "0.8303", 4),
    @(2, "011296", "汇添富优势行业一年定开混合A", "11.39", "86.24", "3.32", "0.3781", 10),
    @(3, "006039", "国富估值优势混合", "6.13", "81.55", "3.16", "0.1937", 3),
    @(4, "000739", "平安新鑫先锋混合A", "5.38", "84.94", "3.26", "0.1754", 6),
    @(5, "010328", "博时荣华灵活配置混合A", "4.31", "77.56", "2.91", "0.1254", 8),
    @(6, "001543", "宝盈新锐灵活配置混合A", "2.20", "91.14", "4.71", "0.1036", 6),
    @(7, "001515", "平安新鑫先锋混合C", "3.15", "84.94", "3.26", "0.1027", 6),
    @(8, "000219", "博时裕益灵活配置混合", "2.12", "87.97", "4.45", "0.0943", 7),
    @(9, "011980", "富兰克林国海匠心精选混合A", "3.74", "85.45", "2.41", "0.0901", 5),
    @(10, "015303", "鹏扬丰融价值先锋一年持有混合A", "3.62", "64.08", "2.47", "0.0894", 5),
    @(11, "005933", "新疆前海联合先进制造灵活配置混合A", "0.95", "92.14", "4.76", "0.0452", 7),
    @(12, "011807", "平安研究精选混合A", "1.12", "89.39", "3.44", "0.0385", 6),
    @(13, "005351", "汇添富行业整合主题混合A", "0.67", "83.14", "5.75", "0.0385", 2),
    @(14, "582003", "东吴配置优化灵活配置混合A", "0.69", "91.71", "4.94", "0.0341", 3),
    @(15, "007581", "宝盈鸿利收益灵活配置混合C", "0.66", "87.42", "5.17", "0.0341", 4),
    @(16, "011707", "东吴配置优化灵活配置混合C", "0.68", "91.71", "4.94", "0.0336", 3),
    @(17, "000066", "诺安鸿鑫混合A", "0.69", "79.09", "4.58", "0.0316", 3),
    @(18, "014151", "国富鑫享价值一年封闭混合A", "2.16", "46.24", "1.32", "0.0285", 4),
    @(19, "011981", "富兰克林国海匠心精选混合C", "0.99", "85.45", "2.41", "0.0239", 5),
    @(20, "004332", "恒生前海沪港深新兴产业精选混合", "0.47", "92.74", "4.71", "0.0221", 6),
    @(21, "011808", "平安研究精选混合C", "0.62", "89.39", "3.44", "0.0213", 6),
    @(22, "007578", "宝盈新锐灵活配置混合C", "0.36", "91.14", "4.71", "0.0170", 6),
    @(23, "014608", "中欧周期景气混合A", "0.24", "90.71", "4.73", "0.0114", 8),
    @(24, "014152", "国富鑫享价值一年封闭混合C", "0.75", "46.24", "1.32", "0.0099", 4),
    @(25, "000788", "前海开源中国成长灵活配置混合", "0.46", "88.66", "2.08", "0.0096", 8),
    @(26, "015304", "鹏扬丰融价值先锋一年持有混合C", "0.36", "64.08", "2.47", "0.0089", 5),
    @(27, "011297", "汇添富优势行业一年定开混合C", "0.15", "86.24", "3.32", "0.0050", 10),
    @(28, "005934", "新疆前海联合先进制造灵活配置混合C", "0.10", "92.14", "4.76", "0.0048", 7),
    @(29, "010329", "博时荣华灵活配置混合C", "0.14", "77.56", "2.91", "0.0041", 8),
    @(30, "014609", "中欧周期景气混合C", "0.06", "90.71", "4.73", "0.0028", 8),
    @(31, "014498", "诺安鸿鑫混合C", "0.01", "79.09", "4.58", "0.0005", 3),
    @(32, "015191", "汇添富行业整合主题混合D", "0.00", "83.14", "5.75", "0", 2),
    @(33, "015190", "汇添富行业整合主题混合C", "0.00", "83.14", "5.75", "0", 2)
)

$rowNum = 2
foreach ($rec in $q3Data) {
    $idx = $rec[0]
    $code = $rec[1]
    $name = $rec[2]
    $scale = $rec[3]
    $position = $rec[4]
    $ratio = $rec[5]
    $marketValue = $rec[6]
    $rank = $rec[7]

    $q3Sheet.Cells.Item($rowNum, 1).Value = $idx

    Set-TextCell $q3Sheet $rowNum 2 $code
    $q3Sheet.Cells.Item($rowNum, 3).Value = $name
    Set-TextCell $q3Sheet $rowNum 4 $scale
    Set-TextCell $q3Sheet $rowNum 5 $position
    Set-TextCell $q3Sheet $rowNum 6 $ratio

    if ($marketValue -eq "0") {
        $q3Sheet.Cells.Item($rowNum, 7).Value = 0
    } else {
        Set-TextCell $q3Sheet $rowNum 7 $marketValue
    }

    $q3Sheet.Cells.Item($rowNum, 8).Value = $rank

    $rowNum = $rowNum + 1
}

# Re-apply the bold/bordered header & index-column formatting
# (direct Value assignment above does not carry cell styles).
$formatSrc.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)

$formatSrc.Range("A2").Copy()
$q3Sheet.Range("A2:A35").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ------------------------------------------------------------
# 2. Prepend the 2022-Q3 summary row on the "总计" sheet, and
#    renumber the existing index column.
# ------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The inserted row inherits stray formatting from neighboring
# rows; reset it to match the source file (plain data cells,
# bold/bordered index cell).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 34
$totalSheet.Cells.Item(2, 4).Value = 4.85

# Renumber column A (index) for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
